$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '56.872.06'
$ws.Range('E2').Value = '  -0.74%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.313.43'
$ws.Range('E3').Value = '  -2.10%  '
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '528.30'
$ws.Range('E5').Value = '  +1.48%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '132.27'
$ws.Range('E6').Value = '  -2.66%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.996'
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').Value = '  -1.28%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.335.12'
$ws.Range('E9').Value = '  -2.06%  '
$ws.Range('E10').Value = '  -1.84%  '
$ws.Range('E11').Value = '  -0.01%  '
$ws.Range('E12').Value = '  -2.85%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.350'
$ws.Range('E13').Value = '  +1.95%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.733.26'
$ws.Range('E14').Value = '  -1.88%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '23.36'
$ws.Range('E15').Value = '  -4.67%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '57.028.80'
$ws.Range('E16').Value = '  -0.50%  '
$ws.Range('E17').Value = '  -2.34%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.342.76'
$ws.Range('E18').Value = '  -1.21%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '334.62'
$ws.Range('E19').Value = '  +1.22%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.40'
$ws.Range('E20').Value = '  -1.69%  '
$ws.Range('E21').Value = '  -1.82%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.79'
$ws.Range('E22').Value = '  +1.06%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.997'
$ws.Range('E23').Value = '  -0.27%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '61.98'
$ws.Range('E24').Value = '  +0.80%  '
$ws.Range('E25').Value = '  +1.17%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.67'
$ws.Range('E26').Value = '  -4.10%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').Value = '  +2.53%  '
$ws.Range('E28').Value = '  -0.81%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '173.51'
$ws.Range('E29').Value = '  +4.08%  '
$ws.Range('E30').Value = '  +1.52%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0₃0722'
$ws.Range('E31').Value = '  -3.44%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.09'
$ws.Range('E32').Value = '  -3.27%  '
$ws.Range('E33').Value = '  -0.93%  '
$ws.Range('E34').Value = '  -0.06%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.992'
$ws.Range('E35').Value = '  -0.21%  '
$ws.Range('E36').Value = '  -4.28%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.926'
$ws.Range('E37').Value = '  +0.60%  '
$ws.Range('E38').Value = '  -1.86%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.06'
$ws.Range('E39').Value = '  +13.89%  '
$ws.Range('B40').Value = 'OKB'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '39.18'
$ws.Range('E40').Value = '  +1.09%  '
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.57'
$ws.Range('E41').Value = '  -1.95%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '148.50'
$ws.Range('E42').Value = '  -1.06%  '
$ws.Range('E43').Value = '  -3.01%  '
$ws.Range('E44').Value = '  -1.80%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '282.66'
$ws.Range('E45').Value = '  -3.31%  '
$ws.Range('E46').Value = '  -0.98%  '
$ws.Range('E47').Value = '  -2.01%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '18.70'
$ws.Range('E48').Value = '  +2.56%  '
$ws.Range('E49').Value = '  -1.70%  '
$ws.Range('E50').Value = '  +5.20%  '
$ws.Range('E51').Value = '  -1.67%  '
